# Working on exergy destruction calculations.
# Insert two new columns (C, D) for Temperature (K) and phi (Carnot factor);
# the old "Note" column (previously C) shifts right to become column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).Insert()

# Headers (write D1 "phi" before C1 "TemperatureK" so the new shared-string
# entries land in the same order as the source edit: 56=phi, 57=TemperatureK)
$ws.Range("D1").Value = "phi"
$ws.Range("C1").Value = "TemperatureK"

# Column C: temperature in Kelvin, either a literal value (sourced from a
# hyperlinked reference) or a simple C->K conversion formula.
$ws.Range("C2").Value = 1144
$ws.Range("C3").Formula = "=600+273.15"
$ws.Range("C4").Formula = "=500+273.15"
$ws.Range("C5").Formula = "=2500+273.15"
$ws.Range("C6").Value = 1573
$ws.Range("C7").Formula = "=25+273.15"
$ws.Range("C8").Formula = "=(1200-32)*5/9 + 273.15"
$ws.Range("C9").Formula = "=200+273.15"
$ws.Range("C10").Formula = "=25+273.15"
$ws.Range("C11").Formula = "=25+273.15"
$ws.Range("C12").Formula = "=25+273.15"
$ws.Range("C13").Formula = "=25+273.15"
$ws.Range("C14").Formula = "=25+273.15"
$ws.Range("C15").Formula = "=25+273.15"
$ws.Range("C16").Formula = "=2500+273.15"
$ws.Range("C17").Formula = "=25+273.15"
$ws.Range("C18").Formula = "=2500+273.15"
$ws.Range("C19").Formula = "=2500+273.15"
$ws.Range("C20").Formula = "=25+273.15"
$ws.Range("C21").Formula = "=2500+273.15"
$ws.Range("C22").Formula = "=2500+273.15"
$ws.Range("C23").Formula = "=2500+273.15"
$ws.Range("C24").Formula = "=2500+273.15"
$ws.Range("C25").Formula = "=2500+273.15"
$ws.Range("C26").Formula = "=2500+273.15"
$ws.Range("C27").Formula = "=2500+273.15"
$ws.Range("C28").Formula = "=2500+273.15"
$ws.Range("C29").Formula = "=2500+273.15"
$ws.Range("C30").Value = 1173
$ws.Range("C31").Value = 1573
$ws.Range("C32").Formula = "=25+273.15"
$ws.Range("C33").Value = 2373
$ws.Range("C34").Formula = "=25+273.15"
$ws.Range("C35").Value = 1144

# Column D: phi = 1 - T0/T, T0 = 298.15 K
for ($r = 2; $r -le 35; $r++) {
    $ws.Range("D$r").Formula = "=1 - 298.15/C$r"
}

# Column widths for the two new columns match column B's.
$ws.Range("C1:D1").ColumnWidth = 11.6

# Inserting the columns copied the neighbouring rows' direct formatting
# (the thin-border "style 3" xf, and the row-20/21 border styles) into the
# freshly created C/D cells; strip that back to Normal everywhere it
# shouldn't be, including the old "Note" column's old border style now
# sitting in column E.
$ws.Range("C20:D21").Style = "Normal"
$ws.Range("A22:E29").Style = "Normal"

# Column insert leaves the three existing hyperlinks pointing at their old
# (now-vacated) C-column cells; rebuild them against the shifted column E.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.advancedenergy.com/getmedia/fffbe6eb-14ab-46b2-8b14-eb04ea9883ff/en-im-steel-coke-oven-application-note.pdf")
$ws.Hyperlinks.Add($ws.Range("E8"), "https://www.terpconnect.umd.edu/~nsw/chbe446/Team0-Sampleb-Project1-FinalReport.pdf See Figure 3 for waste heat temperatures.")
$ws.Hyperlinks.Add($ws.Range("E31"), "https://www.advancedenergy.com/getmedia/fffbe6eb-14ab-46b2-8b14-eb04ea9883ff/en-im-steel-coke-oven-application-note.pdf")
$ws.Range("E6").Style = "Hyperlink"
$ws.Range("E8").Style = "Hyperlink"
$ws.Range("E31").Style = "Hyperlink"

# Match the recorded selection from the source edit.
$ws.Range("C2").Select()
